$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Plan")

# ---- WEEK1 block ----

# Row 4: Install and Setup Linux environment with Centos Minimal
$ws.Range("C4").Value = "13/8/2019"
$ws.Range("F4").Value = ""

# Row 5: Basic Linux (guideline updated) / description updated to "Learn command following..."
$ws.Range("B5").Value = "Basic Linux`n1. Common command`n1. User, Group, Security`n2.  File, Folder, Devices`n3.  Networking: Squid, Webmin, Nginx, Firewall...`n"
$ws.Range("C5").Value = "13/8/2019"
$ws.Range("D5").Value = 16
$ws.Range("E5").Value = "Learn command following:`n1. System information`n2. Device`n3. User & Groups`n4. File & Folder`n5. Compose`n6. Filter content in file `n7. Installation Nginx, Squid"
$ws.Range("F5").Value = ""

# Row 6: How to use Git Hub -- new guideline for Git row
$ws.Range("B6").Value = "How to use Git Hub`n1. Command common`n2.  Manage branches ( merge, rebase)"
$ws.Range("C6").Value = 43693
$ws.Range("D6").Value = 8
$ws.Range("E6").Value = "1. Overview about git`n2. Create new repositories`n3. Command `n4. Rebase and conflict"
$ws.Range("E6").HorizontalAlignment = -4131

# Row 7: How to use Docker
$ws.Range("B7").Value = "How to use Docker"
$ws.Range("C7").Value = 43696
$ws.Range("D7").Value = 8
$ws.Range("E7").Value = "1. How to install docker`n2. How to make image in docker"

# ---- WEEK2 block ----

# Row 13: Learn IAM
$ws.Range("B13").Value = "Learn IAM"
$ws.Range("C13").Value = "20/8/2019"
$ws.Range("E13").Value = "1. What is IAM in AWS?`n2. Why do we need use it?"
$ws.Range("E13").HorizontalAlignment = -4131
$ws.Rows.Item(13).RowHeight = 36.75

# Row 14: Learn S3
$ws.Range("B14").Value = "Learn S3"
$ws.Range("C14").Value = "21/8/2019"
$ws.Range("E14").Value = "1. What is S3 in AWS?`n2. Why do we need use it?"
$ws.Range("E14").HorizontalAlignment = -4131
$ws.Rows.Item(14).RowHeight = 30

# Row 15: Learn EC2, EC2 Auto Scaling
$ws.Range("B15").Value = "Learn EC2, EC2 Auto Scaling"
$ws.Range("C15").Value = "22/9/2019"
$ws.Range("E15").Value = "1. What is EC2 in AWS?`n2. Why do we need use it?"
$ws.Range("E15").HorizontalAlignment = -4131
$ws.Rows.Item(15).RowHeight = 33

# Row 16: Learn  RDS, DynamoDB in AWS
$ws.Range("B16").Value = "Learn  RDS, DynamoDB in AWS"
$ws.Range("C16").Value = "23/8/2019"
$ws.Range("E16").Value = "1. What is RDS, DynamoDB in AWS?`n2. Why do we need use it?"
$ws.Range("E16").HorizontalAlignment = -4131
$ws.Rows.Item(16).RowHeight = 38.25

# Row 17: Learn Virtual Private Cloud(VPC)
$ws.Range("B17").Value = "Learn Virtual Private Cloud(VPC)"
$ws.Range("C17").Value = "26/10/2019"
$ws.Range("E17").Value = "1. What is VPC in AWS?`n2. Why do we need use it?"
$ws.Range("E17").HorizontalAlignment = -4131
$ws.Rows.Item(17).RowHeight = 49.5

# ---- WEEK3 block ----

# Row 22: Learn Elastic Load Balance (ELB)
$ws.Range("B22").Value = "Learn Elastic Load Balance (ELB)"
$ws.Range("C22").Value = "27/8/2019"
$ws.Range("E22").Value = "1. What is ELB in AWS?`n2. Why do we need use it?"

# Row 23: Learn S3  Route 53 DNS
$ws.Range("B23").Value = "Learn S3  Route 53 DNS"
$ws.Range("C23").Value = "28/11/2019"
$ws.Range("E23").Value = "1. What is Route 53 DNS in AWS?`n2. Why do we need use it?"

# Row 24: Learn S3  Monitoring, Metrics, and Logging
$ws.Range("B24").Value = "Learn S3  Monitoring, Metrics, and Logging`n1. CloudTrail`n2. CloudWatch"
$ws.Range("C24").Value = "29/8/2019"
$ws.Range("E24").Value = "1. What is CloudTrail, CloudWatch in AWS?`n2. Why do we need use it?"

# Update the active selection shown when the sheet was last saved
$ws.Range("G6").Select()
